$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45: "Quant" <-> "EnergySwap" swapped with row 46 (coin name/link/price/volume)
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"

# Row 46: "EnergySwap" <-> "Quant" swapped with row 45
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"

$ws.Range("D2").Value = "27.720.16"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.853.11"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "'1.020"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -1.38%  "
$ws.Range("D5").Value = "'319.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("D7").Value = "'0.4369"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'0.8827"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "'21.53"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "1.853.16"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "'6.789"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "'5.490"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "'0.07131"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "'88.31"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.89%  "
$ws.Range("D17").Value = "'1.022"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "'0.000009017"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "'1.017"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "27.728.38"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Value = "'5.268"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "'11.14"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("D24").Value = "2.092.75"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").Value = "'2.035"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.80%  "
$ws.Range("D26").Value = "'156.76"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").Value = "'18.65"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").Value = "'5.436"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("D29").Value = "'1.990"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").Value = "'120.77"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("D31").Value = "'0.09038"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'1.226"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "'0.7702"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "'3.021"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.59%  "
$ws.Range("D35").Value = "'4.559"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "'1.139"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").Value = "'0.01977"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").Value = "'0.05301"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "'2.866"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").Value = "'0.5184"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'6.939"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'8.708"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").Value = "'10.76"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").Value = "'109.92"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "'1.714"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "'0.4727"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("D49").Value = "'1.017"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").Value = "'0.06474"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").Value = "'1.845"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.77%  "
